$wb = $excel.ActiveWorkbook

# -----------------------------------------------------------------------
# "Text Functions" sheet: fill in formulas for columns B:H, rows 2-15
# -----------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Text Functions")

# Row 2 gets the "master" (non-shared) formulas
$ws.Range("B2").Formula = "=LEFT(`$A2, 6)"
$ws.Range("C2").Formula = "=LEFT(`$A2, SEARCH(`"-`", `$A2)-1)"
$ws.Range("D2").Formula = "=MID(`$A2, SEARCH(`"-`", `$A2)+1, 2)"
$ws.Range("E2").Formula = "=IF(ISNUMBER(SEARCH(`"SMALL`",`$A2)),`"Small`",IF(ISNUMBER(SEARCH(`"MEDIUM`",`$A2)),`"Medium`",IF(ISNUMBER(SEARCH(`"LARGE`",`$A2)),`"Large`",IF(ISNUMBER(SEARCH(`"XL`",`$A2)),`"XL`",`"Other`"))))"
$ws.Range("F2").Formula = "=RIGHT(`$A2, LEN(`$A2)-SEARCH(`"_`", `$A2))"
$ws.Range("G2").Formula = "=SUBSTITUTE(`$A2, `"-`", `"|`", 2)"
$ws.Range("H2").Formula = "=MID(`$G2, SEARCH(`"|`", `$G2)+1, 3)"

# Rows 3-15 are filled with the same relative formula so Excel records
# them as a shared formula group (t="shared")
$ws.Range("B3:B15").Formula = "=LEFT(`$A3, 6)"
$ws.Range("C3:C15").Formula = "=LEFT(`$A3, SEARCH(`"-`", `$A3)-1)"
$ws.Range("D3:D15").Formula = "=MID(`$A3, SEARCH(`"-`", `$A3)+1, 2)"
$ws.Range("E3:E15").Formula = "=IF(ISNUMBER(SEARCH(`"SMALL`",`$A3)),`"Small`",IF(ISNUMBER(SEARCH(`"MEDIUM`",`$A3)),`"Medium`",IF(ISNUMBER(SEARCH(`"LARGE`",`$A3)),`"Large`",IF(ISNUMBER(SEARCH(`"XL`",`$A3)),`"XL`",`"Other`"))))"
$ws.Range("F3:F15").Formula = "=RIGHT(`$A3, LEN(`$A3)-SEARCH(`"_`", `$A3))"
$ws.Range("G3:G15").Formula = "=SUBSTITUTE(`$A3, `"-`", `"|`", 2)"
$ws.Range("H3:H15").Formula = "=MID(`$G3, SEARCH(`"|`", `$G3)+1, 3)"

# -----------------------------------------------------------------------
# Move / resize the instructions rectangle shape on "Text Functions"
# -----------------------------------------------------------------------
$shp = $ws.Shapes.Item(1)
$shp.Left = 879.5021484375
$shp.Top = 28.0
$shp.Width = 546.375
$shp.Height = 402.6

# -----------------------------------------------------------------------
# View state: "OFFSET & COUNTA" loses the tabSelected flag (selection
# unchanged) - set this first so it is no longer the active sheet
# -----------------------------------------------------------------------
$wsOffset = $wb.Worksheets.Item("OFFSET & COUNTA")
[void]$wsOffset.Range("D10").Select()

# "Text Functions" becomes the active/selected sheet, with cell E13
# selected (and no frozen/topLeft scroll position) - done last so it
# ends up as the active tab
[void]$ws.Activate()
[void]$ws.Range("E13").Select()
